$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before current column C (x1) to make room for "Rep".
# This shifts old C (x1) -> D and old D (x2) -> E.
$ws.Range("C1").EntireColumn.Insert()

# Header row
$ws.Range("C1").Value = "Rep"
$ws.Range("D1").Value = "x1"
$ws.Range("E1").Value = "x2"

# Set Rep = 1 for the existing 12 rows (rows 2-13)
$ws.Range("C2:C13").Value = 1

# Fix up x1/x2 values for rows 2-13 (D and E columns) per the new dataset
$data1 = @(
    @(6, 8),
    @(4, 6),
    @(8, 12),
    @(2, 6),
    @(3, 8),
    @(-3, 2),
    @(4, 3),
    @(-4, 3),
    @(-3, 2),
    @(-4, -5),
    @(3, -3),
    @(-4, -6)
)

for ($i = 0; $i -lt $data1.Count; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 4).Value = $data1[$i][0]
    $ws.Cells.Item($r, 5).Value = $data1[$i][1]
}

# Append the second replicate block, rows 14-25
$data2 = @(
    @(1, 1, 2, 14, 8),
    @(1, 2, 2, 6, 2),
    @(1, 3, 2, 8, 2),
    @(1, 4, 2, 16, -4),
    @(2, 1, 2, 1, 6),
    @(2, 2, 2, 5, 12),
    @(2, 3, 2, 0, 15),
    @(2, 4, 2, 2, 7),
    @(3, 1, 2, 3, -2),
    @(3, 2, 2, -2, 7),
    @(3, 3, 2, -11, 1),
    @(3, 4, 2, -6, 6)
)

for ($i = 0; $i -lt $data2.Count; $i++) {
    $r = 14 + $i
    $row = $data2[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}

$ws.Range("I5").Select()
